$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A97").Value = "GRT-USD"
